$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$t38 = "Spatial intensity"
$d38 = "A measure of the ratio of events at specific points to a unit of area. Spatial intensity describes the spatially continuous surface of event occurrence. In kernel density estimation, a spatial intensity surface integrates (or sums) to the sample size across a study region."
$t39 = "Spatial density"
$d39 = "A standardized metric of spatial intensity. Related to a probability density function, it is a proportionate indicator of how much of the total events occur in a specific region. In kernel density estimation, the density surface integrates (or sums) to 1 across a study region."
$t40 = "Bandwidth"
$d40 = "A measure of the width or spatial extent of a two-dimensional kernel density estimator. The bandwidth is the key to controlling how much smoothing occurs, with larger bandwidths producing more smooth surfaces, and smaller bandwidths producing less smooth surfaces"
$t41 = "Geographic-weighting"
$t42 = "Kernel density estimator"
$d42 = "A non-parametric way to estimate the probability distribution function of a random variable. In spatial (e.g. 2-d) kernel density estimation, it is a way to describe the spatially continuous variation in the intensity of events (points)."
$d41 = "A method for calculating summary weighted statistics by relying on a kernel density estimator to describe the weights in local summaries."
$t43 = "Homogenous Poisson Point Process"
$d43 = "A spatial statistical assumption that the count of events in an arbitrarily small area is distributed Poisson with mean lambda for all regions"
$t44 = "Inhomogenous Poisson Point Process"
$d44 = "A spatial statistical assumption that the count of events in an arbitrarily small area is distributed Poisson with mean lambda that varies through space as a function of the underlying population at risk. This is true for most spatial epidemiology."
$t45 = "Bandwidth, fixed"
$d45 = "A fixed bandwidth means the width or search radius of the spatial kernel density estimator is constant (fixed) for the full study region"
$t46 = "Bandwidth, adaptive"
$d46 = "An adaptive bandwith means the width or search radius of the spatial kernel density estimator varies or adapts through space, usually to maintain a constant number of points within the window. The result is that in areas with few points there is more smoothing, whereas in areas with many points there is more granularity"

# Module numbers (column A) for rows 38-46
$ws.Cells.Item(38, 1).Value = 6
$ws.Cells.Item(39, 1).Value = 6
$ws.Cells.Item(40, 1).Value = 6
$ws.Cells.Item(41, 1).Value = 6
$ws.Cells.Item(42, 1).Value = 6
$ws.Cells.Item(43, 1).Value = 6
$ws.Cells.Item(44, 1).Value = 6
$ws.Cells.Item(45, 1).Value = 6
$ws.Cells.Item(46, 1).Value = 6

# Row 38
$ws.Cells.Item(38, 2).Value = $t38
$ws.Cells.Item(38, 3).Value = $d38
# Row 39
$ws.Cells.Item(39, 2).Value = $t39
$ws.Cells.Item(39, 3).Value = $d39
# Row 40
$ws.Cells.Item(40, 2).Value = $t40
$ws.Cells.Item(40, 3).Value = $d40
# Row 41/42 terms first, then row 42 def, then row 41 def (matches original authoring order)
$ws.Cells.Item(41, 2).Value = $t41
$ws.Cells.Item(42, 2).Value = $t42
$ws.Cells.Item(42, 3).Value = $d42
$ws.Cells.Item(41, 3).Value = $d41
# Row 43
$ws.Cells.Item(43, 2).Value = $t43
$ws.Cells.Item(43, 3).Value = $d43
# Row 44
$ws.Cells.Item(44, 2).Value = $t44
$ws.Cells.Item(44, 3).Value = $d44
# Row 45
$ws.Cells.Item(45, 2).Value = $t45
$ws.Cells.Item(45, 3).Value = $d45
# Row 46
$ws.Cells.Item(46, 2).Value = $t46
$ws.Cells.Item(46, 3).Value = $d46

# Match the author's final view state: scrolled so row 13 is at top, active cell A47 selected
$excel.Goto($ws.Range("A13"), $true)
$ws.Range("A47").Select()
